$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the end-time for the existing row 69 entry (shifted later -> new value)
$ws.Range("E69").Value = 0.55208333333333337

# Insert a new row above the current blank separator row (row 70),
# shifting the separator / summary rows down by one.
$ws.Rows.Item(70).Insert()

# Fill in the new data row (row 70) with the additional working-hours entry.
$ws.Range("A70").Value = 2014
$ws.Range("B70").Value = 3
$ws.Range("C70").Value = 16
$ws.Range("D70").Value = 0.56944444444444442
$ws.Range("E70").Value = 0.75
$ws.Range("D70").NumberFormat = "hh:mm;@"
$ws.Range("E70").NumberFormat = "hh:mm;@"
$ws.Range("F70").Formula = "=(E70-D70)*24*60"
$ws.Range("F70").NumberFormat = "0"
$ws.Range("G70").Formula = "=F70/60"
$ws.Range("G70").NumberFormat = "0.00"

# Fix up the summary SUM formula which now needs to include the new row 71
# (the blank separator row, now shifted down one from row 70 to row 71).
$ws.Range("F72").Formula = "=SUM(F2:F71)"

$ws.Range("G73").Select()
